$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "implemented"
$ws.Range("E3").Value = "Signal backtest (EOD) implemented: DSL + Top-N momentum ranking + presets + summary results."
$ws.Range("F3").Value = "27/12/2025 02:30"
